# Ticket title changes from "2." to "26. Основные характеристики и идеи
# русской религиозно-идеалистической философии", split across four runs:
# "2", "6", ".", " Основные ... философии" (matching the target OOXML diff,
# which shows the single original run broken into four separate <w:r>
# elements with no run-level formatting).
#
# The engine auto-coalesces adjacent runs that share identical formatting,
# so naive InsertAfter/Range.Text edits on the existing run always collapse
# back into one <w:r>. To get four genuinely separate runs we build the
# pieces in their own paragraphs (which are never auto-merged) and then
# join the paragraphs back together by deleting the paragraph marks
# between them; that leaves the runs as siblings inside a single <w:p>
# without re-coalescing them. Inserting the new paragraphs *before* the
# original one (rather than after) means the final merged paragraph keeps
# the original paragraph's own <w:p> identity/attributes.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Create three empty paragraphs ahead of the original "2." paragraph.
$p1.Range.InsertParagraphBefore()
$p1.Range.InsertParagraphBefore()
$p1.Range.InsertParagraphBefore()

$pA = $d.Paragraphs(1)
$pB = $d.Paragraphs(2)
$pC = $d.Paragraphs(3)
$pD = $d.Paragraphs(4)

$pA.Range.Text = "2"
$pB.Range.Text = "6"
$pC.Range.Text = "."
$pD.Range.Text = " Основные характеристики и идеи русской религиозно-идеалистической философии"

# Join the four paragraphs into one by deleting the paragraph marks that
# separate them. After filling the text above, the mark right after "2"
# always sits at document position 1 (it does not shift as later marks
# are removed), so we can repeat the same deletion three times.
$d.Range(1, 2).Delete()
$d.Range(2, 3).Delete()
$d.Range(3, 4).Delete()
